$d = $word.ActiveDocument

# The 20 "keyword" paragraphs each lose their trailing ", " / "," (and the
# one stray leading space on the second "Mood" paragraph) so only the bare
# word remains.
$finds = @(
  "Browse,", "Search, ", "Events, ", "Parties, ", "Bars, ", "Districts, ", "Music,", "Mood, ", "Social, ", "Drinks,",
  "Browse, ", "Search, ", "Events, ", "Parties, ", "Bars, ", "Districts, ", "Music,", " Mood, ", "Social, ", "Drinks"
)
$replaces = @(
  "Browse", "Search", "Events", "Parties", "Bars", "Districts", "Music", "Mood", "Social", "Drinks",
  "Browse", "Search", "Events", "Parties", "Bars", "Districts", "Music", "Mood", "Social", "Drinks"
)

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $find = $finds[$i - 1]
    $replace = $replaces[$i - 1]
    if ($find -ne $replace) {
        # Scoping Find to the paragraph's own range keeps each replacement
        # local; it also merges the old two-run "Drinks" + "," paragraph
        # (#10) down into the single run the target document expects.
        $p.Range.Find.Execute($find, $false, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
    }
}

# The "_GoBack" bookmark starts out collapsed right after the first
# "Social, " run (paragraph 9) and needs to end up collapsed right after
# "Drinks" in the very last paragraph (20) instead. Re-adding a bookmark
# under an existing name relocates it, so we just need a correctly
# positioned, zero-width Range to hand to Bookmarks.Add.
#
# The host mis-resolves a zero-width Range built exactly at the end of the
# document's content (it snaps back to the very start of the document), so
# we temporarily grow the document by one throwaway character, anchor the
# bookmark just before that character (no longer the document's tail), and
# then delete the throwaway character again.
$endRng = $d.Range($d.Content.End, $d.Content.End)
$endRng.InsertAfter("Z")

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$bmRange = $lastPara.Range.Duplicate
$bmRange.MoveStart(1, 6) | Out-Null   # skip past "Drinks" (6 characters)
$bmRange.MoveEnd(1, -2) | Out-Null    # pull end back before the throwaway "Z" + paragraph mark
$d.Bookmarks.Add("_GoBack", $bmRange)

$cleanupRange = $d.Paragraphs.Item($d.Paragraphs.Count).Range.Duplicate
$cleanupRange.MoveStart(1, 6) | Out-Null
$cleanupRange.Delete() | Out-Null
